$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 28, shifting the existing rows 28-32 down to 29-33.
$ws.Rows("28:28").Insert()

# Populate the newly inserted row 28 with the new weekly price entry.
$ws.Range("A28").Value = 10
$ws.Range("B28").Value = "Vega Modelo de Temuco"
$ws.Range("C28").Value = "La Araucanía"
$ws.Range("D28").Value = 44449
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 300000001
$ws.Range("G28").Value = "Rabanito"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 65
$ws.Range("K28").Value = 7000
$ws.Range("L28").Value = 7000
$ws.Range("M28").Value = 7000
$ws.Range("N28").Value = "$/docena de paquetes"
$ws.Range("O28").Value = "Provincia de Cautín"
$ws.Range("P28").Value = 583
$ws.Range("Q28").Value = 12
$ws.Range("R28").Value = "Hortaliza"
